$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 64 new OHLCV rows (rows 1169-1232) after the existing data ---
# Each entry: (row, timestamp(A), open(B), high(C), low(D), close(E), volume(F))
$newRows = @(
    @(1169, 45534.5, 4.272, 4.314, 4.234, 4.285, 221687.47),
    @(1170, 45534.58333333334, 4.287, 4.302, 4.101, 4.12, 683227.73),
    @(1171, 45534.66666666666, 4.121, 4.221, 4.103, 4.208, 384487.54),
    @(1172, 45534.75, 4.209, 4.276, 4.202, 4.223, 268583.49),
    @(1173, 45534.83333333334, 4.222, 4.256, 4.222, 4.248, 82676.89),
    @(1174, 45534.91666666666, 4.249, 4.291, 4.249, 4.277, 123390.45),
    @(1175, 45535, 4.277, 4.313, 4.276, 4.294, 105834.98),
    @(1176, 45535.08333333334, 4.294, 4.306, 4.277, 4.293, 99833.47),
    @(1177, 45535.16666666666, 4.294, 4.318, 4.293, 4.293, 129625.3),
    @(1178, 45535.25, 4.292, 4.3, 4.288, 4.289, 93956.7),
    @(1179, 45535.33333333334, 4.288, 4.295, 4.259, 4.266, 103219.54),
    @(1180, 45535.41666666666, 4.266, 4.299, 4.265, 4.287, 61812),
    @(1181, 45535.5, 4.286, 4.288, 4.254, 4.257, 114256.87),
    @(1182, 45535.58333333334, 4.256, 4.276, 4.25, 4.269, 77974.22),
    @(1183, 45535.66666666666, 4.27, 4.274, 4.218, 4.24, 113063.24),
    @(1184, 45535.75, 4.24, 4.272, 4.233, 4.241, 93781.17999999999),
    @(1185, 45535.83333333334, 4.24, 4.246, 4.228, 4.245, 53643.8),
    @(1186, 45535.91666666666, 4.261, 4.262, 4.246, 4.258, 30074.54),
    @(1187, 45536, 4.259, 4.259, 4.206, 4.209, 83536.27),
    @(1188, 45536.08333333334, 4.21, 4.225, 4.186, 4.219, 143260.13),
    @(1189, 45536.16666666666, 4.22, 4.225, 4.144, 4.157, 426146.25),
    @(1190, 45536.25, 4.156, 4.203, 4.152, 4.203, 88349.25999999999),
    @(1191, 45536.33333333334, 4.203, 4.206, 4.171, 4.184, 108324.35),
    @(1192, 45536.41666666666, 4.183, 4.193, 4.155, 4.192, 127485.08),
    @(1193, 45536.5, 4.193, 4.199, 4.129, 4.167, 251403.29),
    @(1194, 45536.58333333334, 4.167, 4.195, 4.108, 4.189, 300201.49),
    @(1195, 45536.66666666666, 4.186, 4.201, 4.142, 4.157, 111270.55),
    @(1196, 45536.75, 4.125, 4.221, 4.124, 4.195, 221635.47),
    @(1197, 45536.83333333334, 4.195, 4.217, 4.152, 4.173, 146798.72),
    @(1198, 45536.91666666666, 4.171, 4.175, 4.032, 4.068, 383126.89),
    @(1199, 45537, 4.068, 4.126, 4.067, 4.098, 175623.88),
    @(1200, 45537.08333333334, 4.098, 4.126, 4.085, 4.115, 161703.05),
    @(1201, 45537.16666666666, 4.114, 4.133, 4.066, 4.12, 138993.23),
    @(1202, 45537.25, 4.119, 4.128, 4.051, 4.073, 146381.83),
    @(1203, 45537.33333333334, 4.072, 4.162, 4.072, 4.159, 251619.37),
    @(1204, 45537.41666666666, 4.158, 4.209, 4.157, 4.19, 219143.74),
    @(1205, 45537.5, 4.19, 4.203, 4.177, 4.182, 48454.56),
    @(1206, 45537.58333333334, 4.148, 4.183, 4.122, 4.157, 276257.18),
    @(1207, 45537.66666666666, 4.157, 4.174, 4.135, 4.168, 110399.06),
    @(1208, 45537.75, 4.168, 4.173, 4.138, 4.166, 135390.55),
    @(1209, 45537.83333333334, 4.167, 4.234, 4.163, 4.213, 248261.64),
    @(1210, 45537.91666666666, 4.214, 4.222, 4.199, 4.208, 90136.35000000001),
    @(1211, 45538, 4.209, 4.238, 4.2, 4.233, 95744.91),
    @(1212, 45538.08333333334, 4.233, 4.279, 4.224, 4.259, 175571.92),
    @(1213, 45538.16666666666, 4.259, 4.265, 4.232, 4.264, 115037.54),
    @(1214, 45538.25, 4.262, 4.268, 4.242, 4.265, 111628.62),
    @(1215, 45538.33333333334, 4.264, 4.275, 4.225, 4.245, 188107.17),
    @(1216, 45538.41666666666, 4.246, 4.274, 4.241, 4.265, 103424.19),
    @(1217, 45538.5, 4.264, 4.285, 4.166, 4.193, 265001.02),
    @(1218, 45538.58333333334, 4.192, 4.192, 4.096, 4.108, 391827.15),
    @(1219, 45538.66666666666, 4.107, 4.155, 4.096, 4.125, 196159.22),
    @(1220, 45538.75, 4.125, 4.136, 4.101, 4.106, 144508.99),
    @(1221, 45538.83333333334, 4.107, 4.138, 4.105, 4.13, 130226.3),
    @(1222, 45538.91666666666, 4.131, 4.133, 4.056, 4.062, 127982.32),
    @(1223, 45539, 4.063, 4.097, 3.868, 3.977, 1280889.79),
    @(1224, 45539.08333333334, 3.977, 4.04, 3.965, 4.038, 176739.2),
    @(1225, 45539.16666666666, 4.037, 4.064, 4.032, 4.037, 214690.45),
    @(1226, 45539.25, 4.036, 4.09, 4.036, 4.084, 181957.55),
    @(1227, 45539.33333333334, 4.085, 4.09, 4.046, 4.047, 115643.46),
    @(1228, 45539.41666666666, 4.048, 4.055, 4.012, 4.015, 137661.83),
    @(1229, 45539.5, 4.015, 4.045, 4.005, 4.024, 187789.18),
    @(1230, 45539.58333333334, 4.024, 4.162, 4.004, 4.157, 1020662.84),
    @(1231, 45539.66666666666, 4.156, 4.197, 4.132, 4.133, 279899.8),
    @(1232, 45539.75, 4.134, 4.143, 4.095, 4.109, 206354.24)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# Column A on the new rows needs the same date/time style (border + bold +
# "YYYY-MM-DD HH:MM:SS" number format) used by the existing timestamp column,
# so copy the format from an already-styled cell (A2) onto the new range.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1169:A1232").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
